$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported for "Ají" (Americana (o), Región del
# Maule) dated 2021-12-29. It belongs chronologically right after the
# existing row 36, so insert a fresh row at 37 - this pushes the former
# rows 37-50 down to 38-51 (Excel copies formatting/styles automatically).
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(37, 1).Value = 7
$ws.Cells.Item(37, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(37, 3).Value = "Ñuble"
$ws.Cells.Item(37, 4).Value = 44559
$ws.Cells.Item(37, 5).Value = 16
$ws.Cells.Item(37, 6).Value = 100112021
$ws.Cells.Item(37, 7).Value = "Ají"
$ws.Cells.Item(37, 8).Value = "Americana (o)"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 60
$ws.Cells.Item(37, 11).Value = 18000
$ws.Cells.Item(37, 12).Value = 18500
$ws.Cells.Item(37, 13).Value = 18250
$ws.Cells.Item(37, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(37, 15).Value = "Región del Maule"
$ws.Cells.Item(37, 16).Value = 1217
$ws.Cells.Item(37, 17).Value = 15
$ws.Cells.Item(37, 18).Value = "Hortaliza"

# The source data for the row that ended up at 51 (previously row 50) had
# its "Origen" corrected to "Región del Maule" as part of this update.
$ws.Cells.Item(51, 15).Value = "Región del Maule"
